# Update the "id" and "token" columns (D and C) for rows 2-4 with newly
# generated values, simulating a refreshed configuration / transformer
# timeout run that re-issued ids and JWT tokens for the three test users.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - daniel5f
$ws.Range("D2").Value = "1f3f982e-31a5-4c46-a9b7-ac5c22c5aa70"
$ws.Range("C2").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6ImRhbmllbDVmIiwicGFzc3dvcmQiOiJBejI1Mjg4QCIsImlhdCI6MTcwMTgzNjY0M30.gsKT5xatn2mDcvSFI_f8SPgWPJs_rlZw7v6fVO6CwkU"

# Row 3 - Jorge2525
$ws.Range("D3").Value = "32165355-01d0-4120-8803-d76c3081a121"
$ws.Range("C3").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6IkpvcmdlMjUyNSIsInBhc3N3b3JkIjoiYXNUMzU2NDQ0QCIsImlhdCI6MTcwMTgzNjY0NH0.AC-xAJnICBPfQCgNprypyp987C_CCDKgblk3ek5S73U"

# Row 4 - mario35
$ws.Range("D4").Value = "c7628b39-3ec8-4eba-a281-6ddb65108e69"
$ws.Range("C4").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6Im1hcmlvMzUiLCJwYXNzd29yZCI6Im1BcmlvdXVnQDMiLCJpYXQiOjE3MDE4MzY2NDZ9.yerOfddwtQ2aycT8FZR2pR_kQFSDHhix92gFUbYBroE"
